$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A26").Value = "grouping_steps"
$ws.Range("B26").Value = 3

$ws.Range("A2").Select()
